$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("E6").Value = 64.81999999999999
$ws1.Range("G6").Value = 46.2

$ws1.Range("M10").Value = 24.12

$ws1.Range("M12").Value = 3252.79

$ws1.Range("E13").Value = 64.81999999999999
$ws1.Range("G13").Value = 46.2

$ws1.Range("E21").Value = 129.64
$ws1.Range("H21").Value = 1128.6
$ws1.Range("I21").Value = 432

$ws1.Range("E22").Value = "4 de 20"
$ws1.Range("G22").Value = "4 de 20"
$ws1.Range("H22").Value = "3 de 20"
$ws1.Range("I22").Value = "2 de 20"
$ws1.Range("M22").Value = "6 de 20"

# ---------------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F6").Value = 111.02
$ws2.Range("F10").Value = 4711.21
$ws2.Range("F12").Value = 4381.39
$ws2.Range("F13").Value = 2267.56
$ws2.Range("F21").Value = 1964.99
$ws2.Range("F22").Value = 16249.51

# ---------------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D4").Value = 388.92
$ws3.Range("E4").Value = 124.911046659336
$ws3.Range("F4").Value = 0.7569024926161175

$ws3.Range("D6").Value = 217.42
$ws3.Range("E6").Value = -110.6
$ws3.Range("F6").Value = 2.035386631716907

$ws3.Range("D7").Value = 6013.8
$ws3.Range("E7").Value = -4213.8
$ws3.Range("F7").Value = 3.341

$ws3.Range("D8").Value = 1692
$ws3.Range("E8").Value = -1067
$ws3.Range("F8").Value = 2.7072

$ws3.Range("D16").Value = 5668.87
$ws3.Range("E16").Value = 23863.57
$ws3.Range("F16").Value = 0.1919540004144595

$ws3.Range("D19").Value = 16249.51
$ws3.Range("E19").Value = 34137.68762291769
$ws3.Range("F19").Value = 0.3224928308497397

# Column E width 23 -> 22 (characters). ColumnWidth ~= OOXML width - 1.83
$ws3.Columns.Item(5).ColumnWidth = 21.17
